$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 3 ("Workflow"): "Simulation Parameters" box (shape "Rectangle 5")
# ---------------------------------------------------------------------------
$slide3 = $p.Slides.Item(3)
$paramsShape = $slide3.Shapes.Item("Rectangle 5")
$paramsText = $paramsShape.TextFrame.TextRange

# "Precision" is the 3rd paragraph in this box (after "Simulation
# Parameters" and "Step Size (dt)"). Use Paragraphs(n,1) rather than
# Find(...) to get a range whose InsertAfter/Text assignment performs a
# proper structural (paragraph-mark aware) edit instead of dumping a
# literal CR into the run's text.

# 1) "Precision" -> "Environment effects"
$precisionPara = $paramsText.Paragraphs(3, 1)
$precisionPara.Text = "Environment effects"

# 2) Add a new bullet "…" right after "Environment effects" (formerly
#    "Precision"), pushing the following blank bullets down by one. Doing
#    this as InsertAfter("<CR>…") on the *existing* paragraph's range
#    (rather than InsertBefore on the following empty paragraph) makes the
#    new run inherit the same rPr (solidFill/schemeClr bg1) as the
#    paragraph it follows.
$paramsText2 = $paramsShape.TextFrame.TextRange
$envEffectsPara = $paramsText2.Paragraphs(3, 1)
$envEffectsPara.InsertAfter("`r…")

# ---------------------------------------------------------------------------
# Slide 5 ("Next Steps"): first bullet in the content placeholder
# ---------------------------------------------------------------------------
$slide5 = $p.Slides.Item(5)
$nextStepsShape = $slide5.Shapes.Item("Content Placeholder 2")
$nextStepsText = $nextStepsShape.TextFrame.TextRange

$thrustRun = $nextStepsText.Find("Create functions to more accurately describe thrust and air density", 0)
$thrustRun.Text = "Create functions to more accurately describe changing thrust and changing air density"
